$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header additions (P1, Q1) - match style of existing header cells (bold/border/center)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: swap I/K/M/O values and add P/Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
